$wb = $excel.ActiveWorkbook

# 1. Update the "Date" metadata value (Metadata!B8).
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# 2. On the "Elements" sheet, swap the contents of the two mapping columns:
#      column AK (37) "Mapping: RIM Mapping"
#      column AL (38) "Mapping: Spécification métier vers l'extension ROR AccomodationFamily"
#    for the header row and every data row.
$ws = $wb.Worksheets.Item("Elements")
for ($r = 1; $r -le 6; $r++) {
    $akCell = $ws.Cells.Item($r, 37)
    $alCell = $ws.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# 3. Swap the column widths to match the new (wider) content of column AL vs AK:
#    AK becomes the wide column (~76.53 chars), AL becomes the narrow one (~24.98 chars).
$ws.Columns.Item(37).ColumnWidth = 75.66666666666667
$ws.Columns.Item(38).ColumnWidth = 24.166666666666668
